$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "61.714.85"
$c.NumberFormat = "General"

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +0.06%  "
$c.NumberFormat = "General"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.455.67"
$c.NumberFormat = "General"

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +2.14%  "
$c.NumberFormat = "General"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.NumberFormat = "General"

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.25%  "
$c.NumberFormat = "General"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "576.61"
$c.NumberFormat = "General"

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +0.48%  "
$c.NumberFormat = "General"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "149.04"
$c.NumberFormat = "General"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +8.10%  "
$c.NumberFormat = "General"

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.454.08"
$c.NumberFormat = "General"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +2.05%  "
$c.NumberFormat = "General"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +0.15%  "
$c.NumberFormat = "General"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +0.05%  "
$c.NumberFormat = "General"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.70"
$c.NumberFormat = "General"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +2.91%  "
$c.NumberFormat = "General"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  -0.14%  "
$c.NumberFormat = "General"

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -0.11%  "
$c.NumberFormat = "General"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.045.72"
$c.NumberFormat = "General"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +2.39%  "
$c.NumberFormat = "General"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -1.24%  "
$c.NumberFormat = "General"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "27.38"
$c.NumberFormat = "General"

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +4.83%  "
$c.NumberFormat = "General"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000176"
$c.NumberFormat = "General"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +0.02%  "
$c.NumberFormat = "General"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.458.19"
$c.NumberFormat = "General"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +2.23%  "
$c.NumberFormat = "General"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "61.793.18"
$c.NumberFormat = "General"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +0.24%  "
$c.NumberFormat = "General"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.18"
$c.NumberFormat = "General"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +4.48%  "
$c.NumberFormat = "General"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.13"
$c.NumberFormat = "General"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +1.02%  "
$c.NumberFormat = "General"

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +1.95%  "
$c.NumberFormat = "General"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "383.11"
$c.NumberFormat = "General"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +1.03%  "
$c.NumberFormat = "General"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +1.19%  "
$c.NumberFormat = "General"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.579.41"
$c.NumberFormat = "General"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +1.90%  "
$c.NumberFormat = "General"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.NumberFormat = "General"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -0.14%  "
$c.NumberFormat = "General"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "72.28"
$c.NumberFormat = "General"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +1.17%  "
$c.NumberFormat = "General"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0000125"
$c.NumberFormat = "General"

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -1.00%  "
$c.NumberFormat = "General"

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.177"
$c.NumberFormat = "General"

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +8.67%  "
$c.NumberFormat = "General"

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +2.87%  "
$c.NumberFormat = "General"

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -10.45%  "
$c.NumberFormat = "General"

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.NumberFormat = "General"

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  +0.17%  "
$c.NumberFormat = "General"

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.27"
$c.NumberFormat = "General"

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -0.01%  "
$c.NumberFormat = "General"

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +0.18%  "
$c.NumberFormat = "General"

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c.NumberFormat = "General"

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +1.41%  "
$c.NumberFormat = "General"

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.31"
$c.NumberFormat = "General"

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +1.14%  "
$c.NumberFormat = "General"

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "7.03"
$c.NumberFormat = "General"

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +2.45%  "
$c.NumberFormat = "General"

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  +2.47%  "
$c.NumberFormat = "General"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "166.95"
$c.NumberFormat = "General"

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  +1.42%  "
$c.NumberFormat = "General"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0789"
$c.NumberFormat = "General"

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +2.04%  "
$c.NumberFormat = "General"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "26.26"
$c.NumberFormat = "General"

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  +7.42%  "
$c.NumberFormat = "General"

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +2.82%  "
$c.NumberFormat = "General"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  +0.37%  "
$c.NumberFormat = "General"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -0.28%  "
$c.NumberFormat = "General"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "42.24"
$c.NumberFormat = "General"

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +1.46%  "
$c.NumberFormat = "General"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +1.60%  "
$c.NumberFormat = "General"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -1.09%  "
$c.NumberFormat = "General"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.652.05"
$c.NumberFormat = "General"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +12.54%  "
$c.NumberFormat = "General"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "24.05"
$c.NumberFormat = "General"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +4.67%  "
$c.NumberFormat = "General"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "6.90"
$c.NumberFormat = "General"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +0.24%  "
$c.NumberFormat = "General"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.22"
$c.NumberFormat = "General"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +8.29%  "
$c.NumberFormat = "General"

